$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.04797804163260821
$ws.Range("J2").Value = 0.2267480407458296
$ws.Range("K2").Value = -0.3904699698200783
$ws.Range("L2").Value = 2.722337173295156
